$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textRange = $ws.Range("D2,D3,D5,D6,D8,D12,D13,D14,D15,D17,D18,D19,D20,D21,D23,D26,D27,D28,D29,D31,D32,D33,D37,D38,D39,D40,D42,D43,D47,D48,D49")
$textRange.NumberFormat = "@"

$ws.Range('D2').Value = '67.422.90'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '2.628.82'
$ws.Range('E3').Value = '  -1.68%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '593.98'
$ws.Range('E5').Value = '  -0.64%  '
$ws.Range('D6').Value = '168.11'
$ws.Range('E6').Value = '  +0.98%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '0.533'
$ws.Range('E8').Value = '  -2.41%  '
$ws.Range('E9').Value = '  -1.70%  '
$ws.Range('E10').Value = '  -2.48%  '
$ws.Range('E11').Value = '  +1.24%  '
$ws.Range('D12').Value = '0.363'
$ws.Range('E12').Value = '  +0.67%  '
$ws.Range('D13').Value = '5.22'
$ws.Range('E13').Value = '  -0.10%  '
$ws.Range('D14').Value = '27.61'
$ws.Range('E14').Value = '  -0.84%  '
$ws.Range('D15').Value = '3.109.00'
$ws.Range('E15').Value = '  -1.74%  '
$ws.Range('E16').Value = '  -1.74%  '
$ws.Range('D17').Value = '67.432.52'
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('D18').Value = '2.641.32'
$ws.Range('E18').Value = '  -1.20%  '
$ws.Range('D19').Value = '11.97'
$ws.Range('E19').Value = '  +1.84%  '
$ws.Range('D20').Value = '8.04'
$ws.Range('E20').Value = '  +4.07%  '
$ws.Range('D21').Value = '356.71'
$ws.Range('E21').Value = '  -1.99%  '
$ws.Range('E22').Value = '  -1.82%  '
$ws.Range('D23').Value = '4.66'
$ws.Range('E23').Value = '  -3.56%  '
$ws.Range('E24').Value = '  -4.57%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').Value = '10.31'
$ws.Range('E26').Value = '  +2.19%  '
$ws.Range('D27').Value = '69.63'
$ws.Range('E27').Value = '  -1.95%  '
$ws.Range('D28').Value = '2.760.05'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('E30').Value = '  -1.79%  '
$ws.Range('D31').Value = '546.17'
$ws.Range('E31').Value = '  -2.05%  '
$ws.Range('D32').Value = '7.93'
$ws.Range('E32').Value = '  -1.17%  '
$ws.Range('D33').Value = '1.35'
$ws.Range('E33').Value = '  -3.00%  '
$ws.Range('E34').Value = '  -2.14%  '
$ws.Range('E35').Value = '  +4.77%  '
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('D37').Value = '1.49'
$ws.Range('E37').Value = '  -4.12%  '
$ws.Range('D38').Value = '156.19'
$ws.Range('E38').Value = '  +0.12%  '
$ws.Range('D39').Value = '19.00'
$ws.Range('E39').Value = '  -2.75%  '
$ws.Range('D40').Value = '0.365'
$ws.Range('E40').Value = '  -2.22%  '
$ws.Range('E41').Value = '  -0.75%  '
$ws.Range('D42').Value = '18.29'
$ws.Range('E42').Value = '  +1.94%  '
$ws.Range('D43').Value = '5.21'
$ws.Range('E43').Value = '  -1.99%  '
$ws.Range('E45').Value = '  -3.88%  '
$ws.Range('E46').Value = '  -0.60%  '
$ws.Range('D47').Value = '152.97'
$ws.Range('E47').Value = '  -0.39%  '
$ws.Range('D48').Value = '0.578'
$ws.Range('E48').Value = '  -2.10%  '
$ws.Range('D49').Value = '3.79'
$ws.Range('E49').Value = '  -1.60%  '
$ws.Range('E50').Value = '  -1.50%  '
$ws.Range('E51').Value = '  -1.17%  '

$textRange.ClearFormats()
